# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates (and a few cell additions/removals)
# to the Leve market-profit sheets, matching the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 656.8570999999999
$ws.Range("I11").Value = 656.8570999999999
$ws.Range("K11").Value = 656.8570999999999
$ws.Range("M11").Value = -516.8570999999999

$ws.Range("H53").Value = 502.9091
$ws.Range("I53").Value = 184.5
$ws.Range("J53").Value = 573.6667
$ws.Range("K53").Value = 184.5
$ws.Range("L53").Value = 573.6667
$ws.Range("M53").Value = 452.5
$ws.Range("N53").Value = -1847.6667

$ws.Range("H80").Value = 958
$ws.Range("I80").Value = 1467
$ws.Range("J80").Value = 449
$ws.Range("K80").Value = 4401
$ws.Range("L80").Value = 1347
$ws.Range("M80").Value = -3403.5
$ws.Range("N80").Value = -3343

$ws.Range("H83").Value = 958
$ws.Range("I83").Value = 1467
$ws.Range("J83").Value = 449
$ws.Range("K83").Value = 13203
$ws.Range("L83").Value = 4041
$ws.Range("M83").Value = -8211
$ws.Range("N83").Value = -14025

$ws.Range("H98").Value = 1184.3334
$ws.Range("I98").Value = 1303.875
$ws.Range("K98").Value = 1303.875
$ws.Range("M98").Value = 194.125

$ws.Range("H103").Value = 1166.6666
$ws.Range("I103").Value = 1166.6666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 3499.9998
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -2913.9998
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 1184.3334
$ws.Range("I122").Value = 1303.875
$ws.Range("K122").Value = 3911.625
$ws.Range("M122").Value = -1461.625

$ws.Range("H137").Value = 2222.8667
$ws.Range("I137").Value = 1427.5714
$ws.Range("K137").Value = 4282.7142
$ws.Range("M137").Value = -1732.7142

$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 2101854
$ws.Range("I2").Value = 2451496.2
$ws.Range("K2").Value = 2451496.2
$ws.Range("M2").Value = -2451383.2

$ws.Range("H32").Value = 3328.9824
$ws.Range("I32").Value = 1754.8654
$ws.Range("K32").Value = 1754.8654
$ws.Range("M32").Value = -1467.8654

$ws.Range("H34").Value = 40379.2
$ws.Range("I34").Value = 41900
$ws.Range("K34").Value = 41900
$ws.Range("M34").Value = -41629

$ws.Range("H45").Value = 3898.1
$ws.Range("I45").Value = 4284.3125
$ws.Range("K45").Value = 4284.3125
$ws.Range("M45").Value = -3907.3125

$ws.Range("H74").Value = 66668748
$ws.Range("I74").Value = 83335064
$ws.Range("K74").Value = 83335064
$ws.Range("M74").Value = -83334190

$ws.Range("H77").Value = 66668748
$ws.Range("I77").Value = 83335064
$ws.Range("K77").Value = 416675320
$ws.Range("M77").Value = -416670952

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H97").Value = 478
$ws.Range("J97").Value = 444
$ws.Range("L97").Value = 444
$ws.Range("N97").Value = -1436

$ws.Range("H102").Value = 7697703
$ws.Range("I102").Value = 7697703
$ws.Range("K102").Value = 7697703
$ws.Range("M102").Value = -7696081

$ws.Range("H116").Value = 2101854
$ws.Range("I116").Value = 2451496.2
$ws.Range("K116").Value = 2451496.2
$ws.Range("M116").Value = -2449202.2

$ws.Range("H122").Value = 3985.7585
$ws.Range("J122").Value = 24997
$ws.Range("L122").Value = 74991
$ws.Range("N122").Value = -79891

$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 2101854
$ws.Range("I3").Value = 2451496.2
$ws.Range("K3").Value = 2451496.2
$ws.Range("M3").Value = -2451382.2

$ws.Range("H20").Value = 2416.4736
$ws.Range("I20").Value = 2527
$ws.Range("J20").Value = 2227
$ws.Range("K20").Value = 2527
$ws.Range("L20").Value = 2227
$ws.Range("M20").Value = -2280
$ws.Range("N20").Value = -2721

$ws.Range("H86").Value = 2297.2
$ws.Range("I86").Value = 2624
$ws.Range("J86").Value = 990
$ws.Range("K86").Value = 2624
$ws.Range("L86").Value = 990
$ws.Range("M86").Value = -1501
$ws.Range("N86").Value = -3236

$ws.Range("H89").Value = 2297.2
$ws.Range("I89").Value = 2624
$ws.Range("J89").Value = 990
$ws.Range("K89").Value = 13120
$ws.Range("L89").Value = 4950
$ws.Range("M89").Value = -7504
$ws.Range("N89").Value = -16182

$ws.Range("H94").Value = 1140.3
$ws.Range("I94").Value = 1200.5294
$ws.Range("K94").Value = 1200.5294
$ws.Range("M94").Value = -749.5293999999999

$ws.Range("H99").Value = 2313.96
$ws.Range("I99").Value = 2268.05
$ws.Range("K99").Value = 2268.05
$ws.Range("M99").Value = -770.0500000000002

$ws.Range("H105").Value = 3076
$ws.Range("I105").Value = 2919.0715
$ws.Range("K105").Value = 2919.0715
$ws.Range("M105").Value = -1172.0715

$ws.Range("H107").Value = 144766
$ws.Range("I107").Value = 2227.1667
$ws.Range("K107").Value = 2227.1667
$ws.Range("M107").Value = -307.1667000000002

$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 4004400
$ws.Range("I4").Value = 5500
$ws.Range("K4").Value = 5500
$ws.Range("M4").Value = -5388

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H31").Value = 3611.8928
$ws.Range("I31").Value = 4165.7617
$ws.Range("J31").Value = 1950.2858
$ws.Range("K31").Value = 4165.7617
$ws.Range("L31").Value = 1950.2858
$ws.Range("M31").Value = -3870.7617
$ws.Range("N31").Value = -2540.2858

$ws.Range("H32").Value = 20620.908
$ws.Range("I32").Value = 15654.833
$ws.Range("J32").Value = 26580.2
$ws.Range("K32").Value = 15654.833
$ws.Range("L32").Value = 26580.2
$ws.Range("M32").Value = -15338.833
$ws.Range("N32").Value = -27212.2

$ws.Range("H34").Value = 3611.8928
$ws.Range("I34").Value = 4165.7617
$ws.Range("J34").Value = 1950.2858
$ws.Range("K34").Value = 4165.7617
$ws.Range("L34").Value = 1950.2858
$ws.Range("M34").Value = -3963.7617
$ws.Range("N34").Value = -2354.2858

$ws.Range("H35").Value = 3010
$ws.Range("I35").Value = 1262.5
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 1262.5
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -968.5
$ws.Range("N35").Value = -10588

$ws.Range("H60").Value = 25714.285
$ws.Range("I60").Value = 20000
$ws.Range("J60").Value = 26666.666
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 26666.666
$ws.Range("M60").Value = -19489
$ws.Range("N60").Value = -27688.666

$ws.Range("H86").Value = 8897.323
$ws.Range("I86").Value = 7386.5
$ws.Range("K86").Value = 7386.5
$ws.Range("M86").Value = -6263.5

$ws.Range("H89").Value = 8897.323
$ws.Range("I89").Value = 7386.5
$ws.Range("K89").Value = 36932.5
$ws.Range("M89").Value = -31316.5

$ws.Range("H107").Value = 59756.883
$ws.Range("I107").Value = 385.42856
$ws.Range("J107").Value = 101316.9
$ws.Range("K107").Value = 385.42856
$ws.Range("L107").Value = 101316.9
$ws.Range("M107").Value = 1534.57144
$ws.Range("N107").Value = -105156.9

$ws.Range("H122").Value = 1810.8387
$ws.Range("I122").Value = 1830.7273
$ws.Range("J122").Value = 1762.2222
$ws.Range("K122").Value = 5492.1819
$ws.Range("L122").Value = 5286.6666
$ws.Range("M122").Value = -3042.1819
$ws.Range("N122").Value = -10186.6666

$ws.Range("H134").Value = 27780152
$ws.Range("I134").Value = 27780152
$ws.Range("K134").Value = 83340456
$ws.Range("M134").Value = -83337921

$ws = $wb.Worksheets.Item(5)
$ws.Range("H32").Value = 18812878
$ws.Range("J32").Value = 471.75
$ws.Range("L32").Value = 1415.25
$ws.Range("N32").Value = -1981.25

$ws.Range("H62").Value = 13000.5
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372

$ws.Range("H65").Value = 13000.5
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -141864

$ws.Range("H127").Value = 779.3333
$ws.Range("J127").Value = 779.3333
$ws.Range("L127").Value = 2337.9999
$ws.Range("N127").Value = -12257.9999

$ws.Range("H140").Value = 1782.7142
$ws.Range("I140").Value = 1578.091
$ws.Range("J140").Value = 2533
$ws.Range("K140").Value = 4734.272999999999
$ws.Range("L140").Value = 7599
$ws.Range("M140").Value = 445.7270000000008
$ws.Range("N140").Value = -17959

$ws = $wb.Worksheets.Item(6)
$ws.Range("H5").Value = 100000
$ws.Range("J5").Value = 100000
$ws.Range("L5").Value = 100000
$ws.Range("N5").Value = -100224

$ws.Range("I97").Value = 824.26666
$ws.Range("J97").Value = 1675
$ws.Range("K97").Value = 824.26666
$ws.Range("L97").Value = 1675
$ws.Range("M97").Value = -328.26666
$ws.Range("N97").Value = -2667

$ws.Range("H102").Value = 1666.9032
$ws.Range("I102").Value = 1685.6428
$ws.Range("K102").Value = 1685.6428
$ws.Range("M102").Value = -63.64280000000008

$ws.Range("H113").Value = 65010.375
$ws.Range("I113").Value = 93606.73
$ws.Range("J113").Value = 2098.4
$ws.Range("K113").Value = 93606.73
$ws.Range("L113").Value = 2098.4
$ws.Range("M113").Value = -91436.73
$ws.Range("N113").Value = -6438.4

$ws.Range("H132").Value = 17859590
$ws.Range("I132").Value = 17859590
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 53578770
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -53576240
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 4575
$ws.Range("I7").Value = 4339.4287
$ws.Range("K7").Value = 4339.4287
$ws.Range("M7").Value = -4227.4287

$ws.Range("H32").Value = 704
$ws.Range("I32").Value = 704
$ws.Range("K32").Value = 704
$ws.Range("M32").Value = -387

$ws.Range("H126").Value = 4575
$ws.Range("I126").Value = 4339.4287
$ws.Range("K126").Value = 13018.2861
$ws.Range("M126").Value = -10548.2861

$ws = $wb.Worksheets.Item(8)
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H22").Value = 17499.5
$ws.Range("J22").Value = 17499.5
$ws.Range("L22").Value = 17499.5
$ws.Range("N22").Value = -18085.5

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H107").Value = 358.1875
$ws.Range("I107").Value = 406.83334
$ws.Range("K107").Value = 1220.50002
$ws.Range("M107").Value = 699.4999800000001

$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4450
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3149.6
$ws.Range("I126").Value = 3149.6
$ws.Range("K126").Value = 9448.799999999999
$ws.Range("M126").Value = -6978.799999999999

